# Add another response row for the "ทำอะไร" tag.
#
# The sheet stores (tag, response) pairs as one row per response, with a
# blank separator row between each distinct tag group. The "ชมบอท" tag
# (row 25) sat directly above the "ทำอะไร" group (rows 27-28); inserting
# two new rows right before row 25 makes room for a new blank-separated
# slot and pushes "ชมบอท" and everything below it (ทำอะไร, ชื่ออะไร, หัวข้อ,
# ยินดีครับ, มีคำถาม, ...) down by two rows, matching the target layout
# (old row 25 -> 27, old row 27/28 -> 29/30, etc.) while leaving every
# existing tag/response pair's content untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("25:26").Insert()

# Leave the selection on the newly inserted (still blank) row, matching
# where the editor's cursor ended up after making room for the addition.
$ws.Range("B25").Select()
